# Generate Report for Handoff
#
# Regenerating the localization-status report refreshed the "Latest
# Handoff"/"Latest Handoff Datetime" timestamps for the
# 660d7815-cb70-4420-9c8e-954ba0cf5892 file (row 5 on each sheet),
# since it was (re)handed off:
#   - Overview!D5        : 2016-03-22 08:37:05 -> 2016-03-22 08:37:37
#   - zh-cn!E5 (col E)    : 2016-03-22 08:37:00 -> 2016-03-22 08:37:33
#   - de-de!E5 (col E)    : 2016-03-22 08:37:05 -> 2016-03-22 08:37:37
#
# All other cells/strings are untouched (their shared-string indices
# shift in the saved XML only because two new strings get interned,
# not because their text changes).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-22 08:37:37"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-22 08:37:33"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-22 08:37:37"
